$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H header: cluster_class (same bold/centered header style as G1)
$ws.Range("H1").Value = "cluster_class"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# cluster_class values for rows 2..63 (row 2 = first data row)
$values = @(
    0,      # row 2
    0,      # row 3
    -1,     # row 4
    0,      # row 5
    0,      # row 6
    0,      # row 7
    0,      # row 8
    0,      # row 9
    0,      # row 10
    0,      # row 11
    0,      # row 12
    0,      # row 13
    0,      # row 14
    0,      # row 15
    0,      # row 16
    0,      # row 17
    0,      # row 18
    0,      # row 19
    0,      # row 20
    -1,     # row 21
    0,      # row 22
    0,      # row 23
    0,      # row 24
    0,      # row 25
    0,      # row 26
    0,      # row 27
    0,      # row 28
    0,      # row 29
    0,      # row 30
    0,      # row 31
    0,      # row 32
    0,      # row 33
    0,      # row 34
    0,      # row 35
    0,      # row 36
    0,      # row 37
    0,      # row 38
    0,      # row 39
    0,      # row 40
    0,      # row 41
    0,      # row 42
    0,      # row 43
    0,      # row 44
    0,      # row 45
    0,      # row 46
    0,      # row 47
    0,      # row 48
    0,      # row 49
    0,      # row 50
    0,      # row 51
    0,      # row 52
    0,      # row 53
    0,      # row 54
    0,      # row 55
    0,      # row 56
    0,      # row 57
    0,      # row 58
    0,      # row 59
    0,      # row 60
    0,      # row 61
    0,      # row 62
    -0.03   # row 63
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
